$d = $word.ActiveDocument

# 1) Insert the new paragraphs (empty para, the "Asd sf 23o ınldkg" para with
#    spell-check proofErr markers, the relocated _GoBack bookmark para, and a
#    trailing empty para) at the very start of the document, ahead of the
#    existing "Bugün" paragraph.
$insertRange = $d.Range(0, 0)
$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Asd</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>sf</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> 23o </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>ınldkg</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$insertRange.InsertXML($newBodyXml)

# 2) The _GoBack bookmark used to live at the end of the "Bugün" paragraph;
#    it has now moved to the empty paragraph above, so drop it from "Bugün".
$d.Bookmarks("_GoBack").Delete()
